$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "jezyk" (language) column F is being removed from the student-import
# sheet entirely; everything to its right (grupa, nr tel, email, notatka
# rekrutacyjna) shifts one column to the left.

# 1. Stash the hyperlink cell's formatting (the non-underlined blue style
#    already used for the email column) in a scratch cell that sits in a
#    column left of F, so it is unaffected by the upcoming column delete.
$ws.Range("I2").Copy()
$ws.Range("A100").PasteSpecial(-4122)

# 2. Delete the whole "jezyk" column.
$ws.Range("F1").EntireColumn.Delete()

# 3. The worksheet's hyperlinks collection does not auto-shift with the
#    column delete, so re-create the three mailto: hyperlinks pointing at
#    their new home (email moved from column I to column H).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:olo@gmail.com", "", "", "olo@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:ala@gmail.com", "", "", "ala@gmail.com")
$ws.Hyperlinks.Add($ws.Range("H4"), "mailto:ala@gmail.com", "", "", "ala@gmail.com")

# 4. Hyperlinks.Add applies Excel's built-in underlined "Hyperlink" style;
#    restore the original (non-underlined) look by pasting back the format
#    we stashed in step 1.
$ws.Range("A100").Copy()
$ws.Range("H2:H4").PasteSpecial(-4122)

# 5. Clean up the scratch cell and selection/clipboard state.
$ws.Range("A100").Clear()
$ws.Application.CutCopyMode = $false

# Mirror the author's final selection (column F, row 1 - the column that
# used to hold "jezyk" and now holds "grupa").
$ws.Range("F1").Select()
